$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width change (target stored width = 12; ColumnWidth property differs
# from stored XML width by the 5px/MDW padding offset, so compensate here)
$ws.Columns.Item(1).ColumnWidth = 11.1666666666666

# Data updates: swap D/E columns for most rows, special-case rows 4, 12, 19
$ws.Range("D1").Value = 20
$ws.Range("E1").Value = 41

$ws.Range("D2").Value = 27
$ws.Range("E2").Value = 54

$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 43

$ws.Range("D4").Value = 60
$ws.Range("E4").Value = 80

$ws.Range("D6").Value = 20
$ws.Range("E6").Value = 45

$ws.Range("D7").Value = 15
$ws.Range("E7").Value = 32

$ws.Range("D8").Value = 27
$ws.Range("E8").Value = 54

$ws.Range("D9").Value = 15
$ws.Range("E9").Value = 33

$ws.Range("D10").Value = 15
$ws.Range("E10").Value = 31

# Row 12: remove D12 entirely, E12 stays 31
$ws.Range("D12").ClearContents()

$ws.Range("D14").Value = 15
$ws.Range("E14").Value = 37

$ws.Range("D15").Value = 15
$ws.Range("E15").Value = 37

$ws.Range("D16").Value = 15
$ws.Range("E16").Value = 30

$ws.Range("D17").Value = 15
$ws.Range("E17").Value = 37

$ws.Range("D18").Value = 15
$ws.Range("E18").Value = 39

# Row 19: remove D19 entirely, E19 stays 54
$ws.Range("D19").ClearContents()

# Selection change
$ws.Range("D12").Select()
